$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.995.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.566.53'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.56%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.566.90'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.75%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('E11').Value = '  -3.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.387'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.172.12'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.571.66'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.676.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.79%  '
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '387.61'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('E23').Value = '  +5.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.711.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.53%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  +7.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('E30').Value = '  +4.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.05%  '
$ws.Range('E32').Value = '  +22.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.574.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.02'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.47%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '169.07'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.54'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0809'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.18'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +14.56%  '
$ws.Range('E43').Value = '  +1.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.65'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  +3.50%  '
$ws.Range('E47').Value = '  +6.65%  '
$ws.Range('E48').Value = '  +2.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.496.75'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +12.38%  '
$ws.Range('E50').Value = '  +4.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.31%  '
